$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.976.12'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.557.57'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.93'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.11'
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.778.66'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '1.533.20'
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.99'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '26.975.14'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '0.0₃0707'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.32'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  -3.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.96'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.02'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").Value = '1.422.29'
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("E34").Value = '  +3.80%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +11.22%  '
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.530'
$ws.Range("E39").Value = '  +1.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.808'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.67'
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.32'
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.992'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.78'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").Value = '1.692.03'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.44'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  -0.01%  '

# Restore default (General) formatting on cells where we forced text
$ws.Range("D5").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D48").ClearFormats()
